$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.435.39"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.839.08"
$ws.Range("E3").Value = "  +3.34%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "224.93"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "0.560"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("D8").Value = "32.09"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  +4.16%  "

$ws.Range("D10").Value = "0.0713"
$ws.Range("E10").Value = "  +8.54%  "

$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").Value = "2.102.55"
$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").Value = "1.850.06"
$ws.Range("E13").Value = "  +3.72%  "

$ws.Range("D14").Value = "10.87"
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").Value = "0.648"
$ws.Range("E15").Value = "  +3.63%  "

$ws.Range("D16").Value = "34.454.36"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "4.36"
$ws.Range("E17").Value = "  +3.63%  "

$ws.Range("D18").Value = "69.84"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "251.71"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +8.15%  "

$ws.Range("D21").Value = "11.35"
$ws.Range("E21").Value = "  +9.48%  "

$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "4.29"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").Value = "160.89"
$ws.Range("E25").Value = "  +2.76%  "

$ws.Range("D26").Value = "16.72"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +4.15%  "

$ws.Range("D28").Value = "0.116"
$ws.Range("E28").Value = "  +2.14%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").Value = "0.0537"
$ws.Range("E30").Value = "  +5.05%  "

$ws.Range("D31").Value = "3.82"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").Value = "3.62"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("E34").Value = "  +4.27%  "

$ws.Range("D35").Value = "1.457.26"
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").Value = "0.648"
$ws.Range("E36").Value = "  +3.62%  "

$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("E38").Value = "  +3.01%  "

$ws.Range("D39").Value = "0.969"
$ws.Range("E39").Value = "  +8.89%  "

$ws.Range("D40").Value = "82.39"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("E41").Value = "  -2.40%  "

$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("E43").Value = "  +5.04%  "

$ws.Range("E44").Value = "  +5.28%  "

$ws.Range("D45").Value = "1.998.93"
$ws.Range("E45").Value = "  +3.35%  "

$ws.Range("D48").Value = "106.76"
$ws.Range("E48").Value = "  +8.58%  "

$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("E51").Value = "  +7.69%  "

# Swap rows 46 and 47 (WEMIXToken <-> Kaspa)
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "0.0499"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").Value = "  +0.97%  "

